$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 44201
$ws.Range("B3").Value = 44201
$ws.Range("B4").Value = 44205
$ws.Range("B5").Value = 44210
$ws.Range("B6").Value = 44211

$ws.Range("N5").Select()
